# New crime data collected — weekly refresh of cs-en-us-090pct.xlsx
# Updates the report header (week number + date range) and the week's
# crime-statistics grid (rows 15-31) to the newly collected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text touch-ups (rich-text shared strings). Both runs keep the
# same font/size as the rest of the string, so an in-place character
# replacement is enough — only the digits/date substrings change.
# ---------------------------------------------------------------------

# "Volume 32   Number  30" -> "...  31"
$ws.Range("A8").Characters(21, 2).Text = "31"

# "Report Covering the Week  7/21/2025  Through  7/27/2025"
#                          -> "...7/28/2025  Through  8/3/2025"
$ws.Range("C9").Characters(27, 9).Text = "7/28/2025"
$ws.Range("C9").Characters(47, 9).Text = "8/3/2025"

# ---------------------------------------------------------------------
# Helper: some cells flip from a numeric figure to the sheet's textual
# "no data" placeholders ("0" / "***.*", shared strings already used
# elsewhere on the sheet for the same purpose). Writing those literal
# digits straight into .Value auto-coerces back to a number, so we
# round-trip through a formula (forces text result), freeze it to a
# static value, then re-stamp the donor cell's number format onto it
# so the cell style index is unaffected by the trip.
# ---------------------------------------------------------------------

function Set-TextPlaceholder($cellAddr, $text, $donorStyleCell) {
    $cell = $ws.Range($cellAddr)
    $cell.Formula = '="' + $text + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)   # xlPasteValues: freeze formula -> static text
    $ws.Range($donorStyleCell).Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats: adopt the placeholder style
}

# ---------------------------------------------------------------------
# Row 15 - Rape
# ---------------------------------------------------------------------
Set-TextPlaceholder "D15" "0" "C15"
Set-TextPlaceholder "E15" "***.*" "C15"
Set-TextPlaceholder "F15" "0" "C15"
$ws.Range("H15").Value = -100
$ws.Range("N15").Value = -62.962962962963

# ---------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = -40
$ws.Range("F16").Value = 18
$ws.Range("G16").Value = 10
$ws.Range("H16").Value = 80
$ws.Range("I16").Value = 118
$ws.Range("J16").Value = 111
$ws.Range("K16").Value = 6.306306306306
$ws.Range("L16").Value = 5.357142857142
$ws.Range("M16").Value = -42.995169082125
$ws.Range("N16").Value = -84.308510638297

# ---------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 60
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -4
$ws.Range("I17").Value = 168
$ws.Range("J17").Value = 183
$ws.Range("K17").Value = -8.196721311475
$ws.Range("L17").Value = -9.677419354838
$ws.Range("M17").Value = 46.086956521739
$ws.Range("N17").Value = -47.663551401869

# ---------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------
$ws.Range("C18").Value = 3
$ws.Range("D18").Value = 7
$ws.Range("E18").Value = -57.142857142857
$ws.Range("F18").Value = 17
$ws.Range("G18").Value = 20
$ws.Range("H18").Value = -15
$ws.Range("I18").Value = 153
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = -23.5
$ws.Range("L18").Value = 4.794520547945
$ws.Range("M18").Value = -38.8
$ws.Range("N18").Value = -80.583756345177

# ---------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------
$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = 18.181818181818
$ws.Range("F19").Value = 52
$ws.Range("G19").Value = 60
$ws.Range("H19").Value = -13.333333333333
$ws.Range("I19").Value = 407
$ws.Range("J19").Value = 433
$ws.Range("K19").Value = -6.004618937644
$ws.Range("L19").Value = -9.354120267260
$ws.Range("M19").Value = 48
$ws.Range("N19").Value = 32.573289902280

# ---------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------
Set-TextPlaceholder "C20" "0" "C15"
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 14
$ws.Range("G20").Value = 7
$ws.Range("H20").Value = 100
$ws.Range("J20").Value = 69
$ws.Range("K20").Value = 20.289855072463
$ws.Range("L20").Value = -12.631578947368
$ws.Range("M20").Value = -6.741573033707
$ws.Range("N20").Value = -83.945841392649

# ---------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 30
$ws.Range("E21").Value = -10
$ws.Range("F21").Value = 125
$ws.Range("G21").Value = 123
$ws.Range("H21").Value = 1.626016260162
$ws.Range("I21").Value = 941
$ws.Range("J21").Value = 1011
$ws.Range("K21").Value = -6.923837784371
$ws.Range("L21").Value = -5.9
$ws.Range("M21").Value = 0
$ws.Range("N21").Value = -65.442526625045

# ---------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------
Set-TextPlaceholder "F22" "0" "C15"
$ws.Range("H22").Value = -100

# ---------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------
$ws.Range("F23").Value = 10
$ws.Range("G23").Value = 4
$ws.Range("H23").Value = 150
$ws.Range("I23").Value = 94
$ws.Range("J23").Value = 96
$ws.Range("K23").Value = -2.083333333333
$ws.Range("L23").Value = -23.577235772357
$ws.Range("M23").Value = 13.253012048192

# ---------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------
$ws.Range("C24").Value = 25
$ws.Range("D24").Value = 30
$ws.Range("E24").Value = -16.666666666666
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 92
$ws.Range("H24").Value = 4.347826086956
$ws.Range("I24").Value = 678
$ws.Range("J24").Value = 667
$ws.Range("K24").Value = 1.649175412293
$ws.Range("L24").Value = 7.790143084260
$ws.Range("M24").Value = -5.833333333333

# ---------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 24
$ws.Range("H25").Value = 20.833333333333
$ws.Range("I25").Value = 186
$ws.Range("J25").Value = 223
$ws.Range("K25").Value = -16.591928251121
$ws.Range("L25").Value = 135.443037974684

# ---------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------
$ws.Range("C26").Value = 12
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 300
$ws.Range("F26").Value = 51
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 70
$ws.Range("I26").Value = 238
$ws.Range("J26").Value = 321
$ws.Range("K26").Value = -25.856697819314
$ws.Range("L26").Value = -16.491228070175
$ws.Range("M26").Value = -19.865319865319

# ---------------------------------------------------------------------
# Row 27 - UCR Rape*
# ---------------------------------------------------------------------
Set-TextPlaceholder "D27" "0" "C15"
Set-TextPlaceholder "E27" "***.*" "C15"
$ws.Range("F27").Value = 1
$ws.Range("H27").Value = 0

# ---------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------
$ws.Range("C28").Value = 3
Set-TextPlaceholder "D28" "0" "C15"
Set-TextPlaceholder "E28" "***.*" "C15"
$ws.Range("G28").Value = 6
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 37
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = 2.777777777777

# ---------------------------------------------------------------------
# Row 29 - Shooting Vic.
# ---------------------------------------------------------------------
$ws.Range("L29").Value = -75
$ws.Range("N29").Value = -95.555555555555

# ---------------------------------------------------------------------
# Row 30 - Shooting Inc.
# ---------------------------------------------------------------------
$ws.Range("L30").Value = -60
$ws.Range("N30").Value = -95.121951219512

# ---------------------------------------------------------------------
# Row 31 - Hate Crimes
# ---------------------------------------------------------------------
Set-TextPlaceholder "D31" "0" "C15"
Set-TextPlaceholder "E31" "***.*" "C15"
$ws.Range("F31").Value = 2
$ws.Range("H31").Value = 100
$ws.Range("I31").Value = 15
$ws.Range("K31").Value = -6.25
$ws.Range("L31").Value = 66.666666666666
